$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 6
$ws.Range("H8").Value = 6

$ws.Range("F17").Value = 42
$ws.Range("H17").Value = 42

$ws.Range("E18").Value = 97

$ws.Range("E19").Value = 42

$ws.Range("E28").Value = 13

$ws.Range("E33").Value = 31
$ws.Range("F33").Value = 9
$ws.Range("H33").Value = 9

$ws.Range("E34").Value = 14

$ws.Range("E37").Value = 42

$ws.Range("E38").Value = 58

$ws.Range("E40").Value = 16
$ws.Range("F40").Value = 8
$ws.Range("H40").Value = 8

$ws.Range("E41").Value = 30

$ws.Range("E47").Value = 50
$ws.Range("F47").Value = 31
$ws.Range("H47").Value = 31

$ws.Range("E50").Value = 22

$ws.Range("E57").Value = 10

$ws.Range("E72").Value = 34

$ws.Range("E79").Value = 27

$ws.Range("E81").Value = 12

$ws.Range("E83").Value = 10
$ws.Range("F83").Value = 2
$ws.Range("H83").Value = 2

$ws.Range("E88").Value = 18
